# Auto-generated edit script: apply numeric corrections to Leve profit sheets
# per commit "chore: update Sheets via scheduled runner"
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 180
$ws.Range("I18").Value = 180
$ws.Range("K18").Value = 180
$ws.Range("M18").Value = 104
$ws.Range("H40").Value = 75914.80499999999
$ws.Range("J40").Value = 3374.963
$ws.Range("L40").Value = 3374.963
$ws.Range("N40").Value = -3724.963
$ws.Range("H57").Value = 69420
$ws.Range("J57").Value = 69420
$ws.Range("L57").Value = 208260
$ws.Range("N57").Value = -209258
$ws.Range("H76").Value = 5880.1
$ws.Range("J76").Value = 7256.2
$ws.Range("L76").Value = 7256.2
$ws.Range("N76").Value = -7886.2
$ws.Range("H79").Value = 5880.1
$ws.Range("J79").Value = 7256.2
$ws.Range("L79").Value = 7256.2
$ws.Range("N79").Value = -9440.200000000001
$ws.Range("H98").Value = 2005.3549
$ws.Range("I98").Value = 2005.3549
$ws.Range("K98").Value = 2005.3549
$ws.Range("M98").Value = -507.3549
$ws.Range("H112").Value = 1032206.2
$ws.Range("J112").Value = 1130347.4
$ws.Range("L112").Value = 3391042.2
$ws.Range("N112").Value = -3393258.2
$ws.Range("H122").Value = 2005.3549
$ws.Range("I122").Value = 2005.3549
$ws.Range("K122").Value = 6016.0647
$ws.Range("M122").Value = -3566.0647
$ws.Range("H138").Value = 3236.5974
$ws.Range("I138").Value = 2513.8057
$ws.Range("J138").Value = 3871.244
$ws.Range("K138").Value = 7541.4171
$ws.Range("L138").Value = 11613.732
$ws.Range("M138").Value = -2401.4171
$ws.Range("N138").Value = -21893.732
$ws.Range("H141").Value = 8029.3335
$ws.Range("I141").Value = 8029.3335
$ws.Range("K141").Value = 24088.0005
$ws.Range("M141").Value = -18908.0005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H58").Value = 79
$ws.Range("I58").Value = 79
$ws.Range("K58").Value = 79
$ws.Range("M58").Value = 351
$ws.Range("H74").Value = 2067.5588
$ws.Range("J74").Value = 3199.6667
$ws.Range("L74").Value = 3199.6667
$ws.Range("N74").Value = -4947.6667
$ws.Range("H77").Value = 2067.5588
$ws.Range("J77").Value = 3199.6667
$ws.Range("L77").Value = 15998.3335
$ws.Range("N77").Value = -24734.3335
$ws.Range("H113").Value = 60989
$ws.Range("J113").Value = 60989
$ws.Range("L113").Value = 60989
$ws.Range("N113").Value = -69667
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 2405.6667
$ws.Range("J29").Value = 5018
$ws.Range("L29").Value = 5018
$ws.Range("N29").Value = -5596
$ws.Range("H107").Value = 858.7895
$ws.Range("I107").Value = 747
$ws.Range("J107").Value = 1278
$ws.Range("K107").Value = 747
$ws.Range("L107").Value = 1278
$ws.Range("M107").Value = 1173
$ws.Range("N107").Value = -5118

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 277.66666
$ws.Range("J7").Value = 364.57144
$ws.Range("L7").Value = 364.57144
$ws.Range("N7").Value = -590.5714399999999
$ws.Range("H22").Value = 383.7143
$ws.Range("J22").Value = 248
$ws.Range("L22").Value = 248
$ws.Range("N22").Value = -948
$ws.Range("H31").Value = 3915.9744
$ws.Range("I31").Value = 2546.2964
$ws.Range("J31").Value = 6997.75
$ws.Range("K31").Value = 2546.2964
$ws.Range("L31").Value = 6997.75
$ws.Range("M31").Value = -2251.2964
$ws.Range("N31").Value = -7587.75
$ws.Range("H34").Value = 3915.9744
$ws.Range("I34").Value = 2546.2964
$ws.Range("J34").Value = 6997.75
$ws.Range("K34").Value = 2546.2964
$ws.Range("L34").Value = 6997.75
$ws.Range("M34").Value = -2344.2964
$ws.Range("N34").Value = -7401.75
$ws.Range("H134").Value = 1889.1428
$ws.Range("I134").Value = 1734.7407
$ws.Range("K134").Value = 5204.2221
$ws.Range("M134").Value = -2669.2221

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 963.4666999999999
$ws.Range("I5").Value = 732.3889
$ws.Range("J5").Value = 1310.0834
$ws.Range("K5").Value = 2197.1667
$ws.Range("L5").Value = 3930.2502
$ws.Range("M5").Value = -2085.1667
$ws.Range("N5").Value = -4154.2502
$ws.Range("H38").Value = 188.7931
$ws.Range("I38").Value = 195.53847
$ws.Range("J38").Value = 183.3125
$ws.Range("K38").Value = 586.61541
$ws.Range("L38").Value = 549.9375
$ws.Range("M38").Value = -239.61541
$ws.Range("N38").Value = -1243.9375
$ws.Range("H55").Value = 3754.2
$ws.Range("I55").Value = 2004
$ws.Range("K55").Value = 6012
$ws.Range("M55").Value = -5835
$ws.Range("H104").Value = 6922.222
$ws.Range("J104").Value = 7721.875
$ws.Range("L104").Value = 23165.625
$ws.Range("N104").Value = -28407.625
$ws.Range("H123").Value = 4675
$ws.Range("J123").Value = 4675
$ws.Range("L123").Value = 14025
$ws.Range("N123").Value = -18925
$ws.Range("H132").Value = 2951
$ws.Range("I132").Value = 2962.4
$ws.Range("K132").Value = 26661.6
$ws.Range("M132").Value = -24131.6
$ws.Range("H135").Value = 963.4666999999999
$ws.Range("I135").Value = 732.3889
$ws.Range("J135").Value = 1310.0834
$ws.Range("K135").Value = 6591.5001
$ws.Range("L135").Value = 11790.7506
$ws.Range("M135").Value = -4056.5001
$ws.Range("N135").Value = -16860.7506

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1517.909
$ws.Range("I16").Value = 1455.875
$ws.Range("J16").Value = 1683.3334
$ws.Range("K16").Value = 1455.875
$ws.Range("L16").Value = 1683.3334
$ws.Range("M16").Value = -1285.875
$ws.Range("N16").Value = -2023.3334
$ws.Range("H121").Value = 64660
$ws.Range("J121").Value = 64660
$ws.Range("L121").Value = 64660
$ws.Range("N121").Value = -68154
$ws.Range("H132").Value = 4047.7454
$ws.Range("I132").Value = 3925
$ws.Range("K132").Value = 11775
$ws.Range("M132").Value = -9245
$ws.Range("H136").Value = 4946.205
$ws.Range("I136").Value = 4556.7095
$ws.Range("K136").Value = 13670.1285
$ws.Range("M136").Value = -11120.1285

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 28000
$ws.Range("I9").Value = 28000
$ws.Range("K9").Value = 28000
$ws.Range("M9").Value = -27860
$ws.Range("H20").Value = 18258
$ws.Range("I20").Value = 6505
$ws.Range("K20").Value = 6505
$ws.Range("M20").Value = -6265
$ws.Range("H45").Value = 20498.334
$ws.Range("I45").Value = 4250
$ws.Range("J45").Value = 23748
$ws.Range("K45").Value = 4250
$ws.Range("L45").Value = 23748
$ws.Range("M45").Value = -3759
$ws.Range("N45").Value = -24730
$ws.Range("H95").Value = 39085.5
$ws.Range("J95").Value = 39085.5
$ws.Range("L95").Value = 39085.5
$ws.Range("N95").Value = -44577.5
